$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
# Row 17
$ws.Range("H17").Value = 1237.88
$ws.Range("J17").Value = 1237.88
$ws.Range("L17").Value = 3713.64
$ws.Range("N17").Value = -4049.64
# Row 101
$ws.Range("H101").Value = 16667425
$ws.Range("I101").Value = 50000276
$ws.Range("J101").Value = 998.25
$ws.Range("K101").Value = 150000828
$ws.Range("L101").Value = 2994.75
$ws.Range("M101").Value = -149999206
$ws.Range("N101").Value = -6238.75
# Row 113
$ws.Range("H113").Value = 3875.6
$ws.Range("I113").Value = 3289.5
$ws.Range("K113").Value = 3289.5
$ws.Range("M113").Value = -35.5
# Row 129
$ws.Range("H129").Value = 2712.1052
$ws.Range("J129").Value = 2848.889
$ws.Range("L129").Value = 8546.667000000001
$ws.Range("N129").Value = -18546.667
# Row 141
$ws.Range("H141").Value = 3778.1428
$ws.Range("I141").Value = 3812.5
$ws.Range("J141").Value = 3732.3333
$ws.Range("K141").Value = 11437.5
$ws.Range("L141").Value = 11196.9999
$ws.Range("M141").Value = -6257.5
$ws.Range("N141").Value = -21556.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2316.5
$ws.Range("I61").Value = 2316.5
$ws.Range("K61").Value = 2316.5
$ws.Range("M61").Value = -2104.5
# Row 74
$ws.Range("H74").Value = 750
$ws.Range("I74").Value = 750
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 750
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 750
$ws.Range("I77").Value = 750
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 3750
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()
# Row 122
$ws.Range("H122").Value = 2737.25
$ws.Range("I122").Value = 2049.6667
$ws.Range("K122").Value = 6149.000100000001
$ws.Range("M122").Value = -3699.000100000001
# Row 132
$ws.Range("H132").Value = 2477.2
$ws.Range("I132").Value = 2504.1428
$ws.Range("K132").Value = 7512.428400000001
$ws.Range("M132").Value = -4982.428400000001
# Row 136
$ws.Range("H136").Value = 2316.5
$ws.Range("I136").Value = 2316.5
$ws.Range("K136").Value = 6949.5
$ws.Range("M136").Value = -4399.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2368.5
$ws.Range("I134").Value = 2428.2856
$ws.Range("J134").Value = 1950
$ws.Range("K134").Value = 7284.8568
$ws.Range("L134").Value = 5850
$ws.Range("M134").Value = -4749.8568
$ws.Range("N134").Value = -10920

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 309.625
$ws.Range("I22").Value = 288
$ws.Range("J22").Value = 374.5
$ws.Range("K22").Value = 288
$ws.Range("L22").Value = 374.5
$ws.Range("M22").Value = 62
$ws.Range("N22").Value = -1074.5
# Row 23
$ws.Range("H23").Value = 129999.5
$ws.Range("I23").Value = 129999.5
$ws.Range("K23").Value = 129999.5
$ws.Range("M23").Value = -129759.5
# Row 27
$ws.Range("H27").Value = 129999.5
$ws.Range("I27").Value = 129999.5
$ws.Range("K27").Value = 129999.5
$ws.Range("M27").Value = -129807.5
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").ClearContents()
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# Row 92
$ws.Range("H92").Value = 424
$ws.Range("I92").Value = 220.4
$ws.Range("J92").Value = 537.1111
$ws.Range("K92").Value = 661.2
$ws.Range("L92").Value = 1611.3333
$ws.Range("M92").Value = 586.8
$ws.Range("N92").Value = -4107.3333
# Row 107
$ws.Range("H107").Value = 643.1429000000001
$ws.Range("I107").Value = 863
$ws.Range("J107").Value = 350
$ws.Range("K107").Value = 2589
$ws.Range("L107").Value = 1050
$ws.Range("M107").Value = -669
$ws.Range("N107").Value = -4890
# Row 137
$ws.Range("H137").Value = 2099.111
$ws.Range("I137").Value = 1199.5
$ws.Range("J137").Value = 2356.1428
$ws.Range("K137").Value = 3598.5
$ws.Range("L137").Value = 7068.428400000001
$ws.Range("M137").Value = 1501.5
$ws.Range("N137").Value = -17268.4284

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 1833733.4
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 102
$ws.Range("H102").Value = 2648.2222
$ws.Range("I102").Value = 2648.2222
$ws.Range("K102").Value = 2648.2222
$ws.Range("M102").Value = -1026.2222
# Row 122
$ws.Range("H122").Value = 1732.3334
$ws.Range("I122").Value = 1649
$ws.Range("J122").Value = 1899
$ws.Range("K122").Value = 4947
$ws.Range("L122").Value = 5697
$ws.Range("M122").Value = -2497
$ws.Range("N122").Value = -10597
# Row 132
$ws.Range("H132").Value = 4497
$ws.Range("I132").Value = 4497
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13491
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3849.3
$ws.Range("I7").Value = 1927.8572
$ws.Range("K7").Value = 1927.8572
$ws.Range("M7").Value = -1815.8572
# Row 16
$ws.Range("H16").Value = 568.75
$ws.Range("I16").Value = 568.75
$ws.Range("K16").Value = 568.75
$ws.Range("M16").Value = -398.75
# Row 55
$ws.Range("H55").Value = 1652.8096
$ws.Range("I55").Value = 1322.7858
$ws.Range("J55").Value = 2312.8572
$ws.Range("K55").Value = 1322.7858
$ws.Range("L55").Value = 2312.8572
$ws.Range("M55").Value = -1149.7858
$ws.Range("N55").Value = -2658.8572
# Row 122
$ws.Range("H122").Value = 6665.4165
$ws.Range("I122").Value = 5181.727
$ws.Range("J122").Value = 7920.846
$ws.Range("K122").Value = 15545.181
$ws.Range("L122").Value = 23762.538
$ws.Range("M122").Value = -13095.181
$ws.Range("N122").Value = -28662.538
# Row 126
$ws.Range("H126").Value = 3849.3
$ws.Range("I126").Value = 1927.8572
$ws.Range("K126").Value = 5783.571599999999
$ws.Range("M126").Value = -3313.571599999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2061.9167
$ws.Range("I122").Value = 2027
$ws.Range("J122").Value = 2166.6667
$ws.Range("K122").Value = 6081
$ws.Range("L122").Value = 6500.000100000001
$ws.Range("M122").Value = -3631
$ws.Range("N122").Value = -11400.0001
# Row 126
$ws.Range("H126").Value = 5045.486
$ws.Range("I126").Value = 4769.727
$ws.Range("J126").Value = 5512.154
$ws.Range("K126").Value = 14309.181
$ws.Range("L126").Value = 16536.462
$ws.Range("M126").Value = -11839.181
$ws.Range("N126").Value = -21476.462
